$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 33, shifting existing rows 33-35 down to 34-36
$ws.Rows.Item(33).Insert()

# Populate the new row 33 with the "Robin" vessel entry
$ws.Range("A33").Value = "Robin"
$ws.Range("B33").Value = "Edison Chouest Offshore"
$ws.Range("C33").Value = 280
$ws.Range("D33").Value = "OSV"

# Update the view state to reflect the final selection/scroll position
$excel.ActiveWindow.ScrollRow = 3
$ws.Range("A34").Select()
